$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the empty, unused Title placeholder shape from the slide.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    $isTitlePlaceholder = $false
    try {
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 1) {
            $isTitlePlaceholder = $true
        }
    } catch {
        $isTitlePlaceholder = $false
    }
    if ($isTitlePlaceholder) {
        $sh.Delete()
    }
}

Write-Output "Shapes remaining: $($s.Shapes.Count)"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    Write-Output "$i : $($sh.Name)"
}
